# Add a "Total clocks" summary row (row 25) to each of the four
# philosopher-algorithm worksheets: per-column SUM(..2:..21) plus a
# combined grand-total SUM across all philosopher columns, labeled with
# a new "Total clocks" shared string.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Algorithm 2, 4 Philosophers" (4 data columns A:D) -----------
$ws1 = $wb.Worksheets.Item("Algorithm 2, 4 Philosophers")
$ws1.Activate()
$ws1.Range("A25").Formula = "=SUM(A2:A21)"
$ws1.Range("B25:D25").Formula = "=SUM(B2:B21)"
$ws1.Range("E25").Formula = "=SUM(A2:D21)"
$ws1.Range("F25").Value = "Total clocks"
$ws1.Range("A25").Select()

# --- Sheet 2: "Algorithm 2, 5 Philosophers" (5 data columns A:E) ----------
$ws2 = $wb.Worksheets.Item("Algorithm 2, 5 Philosophers")
$ws2.Activate()
$ws2.Range("A25").Formula = "=SUM(A2:A21)"
$ws2.Range("B25:E25").Formula = "=SUM(B2:B21)"
$ws2.Range("F25").Formula = "=SUM(A2:E21)"
$ws2.Range("G25").Value = "Total clocks"
$ws2.Range("A25").Select()

# --- Sheet 3: "Algorithm 3, 4 Philosophers" (4 data columns A:D) ----------
$ws3 = $wb.Worksheets.Item("Algorithm 3, 4 Philosophers")
$ws3.Activate()
$ws3.Range("A25").Formula = "=SUM(A2:A21)"
$ws3.Range("B25:D25").Formula = "=SUM(B2:B21)"
$ws3.Range("E25").Formula = "=SUM(A2:D21)"
$ws3.Range("F25").Value = "Total clocks"
$ws3.Range("A25:D25").Select()

# --- Sheet 4: "Algorithm 3, 5 Philosophers" (5 data columns A:E) ----------
$ws4 = $wb.Worksheets.Item("Algorithm 3, 5 Philosophers")
$ws4.Activate()
$ws4.Range("A25").Formula = "=SUM(A2:A21)"
$ws4.Range("B25:E25").Formula = "=SUM(B2:B21)"
$ws4.Range("F25").Formula = "=SUM(A2:E21)"
$ws4.Range("G25").Value = "Total clocks"
$ws4.Range("A25:E25").Select()

# End with "Algorithm 2, 5 Philosophers" as the active tab / selection,
# matching where the author ended their editing session.
$ws2.Activate()
$ws2.Range("A25").Select()
